# Add a new column K "intervention_type" with header style matching the
# other header cells (A1:J1), and fill in values for rows 2-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1 - same style as the other header cells (bold, centered, bordered)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("K1").Style = $ws.Range("A1").Style

# Data values for K2:K14
$values = @{
    2  = "BEHAVIORAL"
    3  = "BEHAVIORAL"
    4  = "BEHAVIORAL"
    5  = "BEHAVIORAL"
    6  = "DEVICE"
    7  = "DEVICE"
    8  = "OTHER"
    9  = "OTHER"
    10 = "OTHER"
    11 = "BEHAVIORAL"
    12 = "OTHER"
    13 = "BEHAVIORAL"
    14 = "OTHER"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 11).Value = $values[$row]
}
